# Auto-generated Excel COM-interop script
# Applies targeted cell value updates to the Leviathan_Profits style workbook
# (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1918351.9
$ws.Range("I33").Value = 3448553.5
$ws.Range("K33").Value = 3448553.5
$ws.Range("M33").Value = -3448324.5
$ws.Range("H43").Value = 10699.75
$ws.Range("J43").Value = 11266.333
$ws.Range("L43").Value = 11266.333
$ws.Range("N43").Value = -11404.333
$ws.Range("H51").Value = 4391385.5
$ws.Range("J51").Value = 4635240.5
$ws.Range("L51").Value = 4635240.5
$ws.Range("N51").Value = -4636208.5
$ws.Range("H125").Value = 36152.25
$ws.Range("I125").Value = 53511.375
$ws.Range("J125").Value = 1434
$ws.Range("K125").Value = 481602.375
$ws.Range("L125").Value = 12906
$ws.Range("M125").Value = -479142.375
$ws.Range("N125").Value = -17826
$ws.Range("H137").Value = 5661.029
$ws.Range("I137").Value = 2260.5334
$ws.Range("K137").Value = 6781.600199999999
$ws.Range("M137").Value = -4231.600199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3400.8333
$ws.Range("I2").Value = 2882
$ws.Range("J2").Value = 5995
$ws.Range("K2").Value = 2882
$ws.Range("L2").Value = 5995
$ws.Range("M2").Value = -2769
$ws.Range("N2").Value = -6221
$ws.Range("H34").Value = 16859
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H74").Value = 2324.276
$ws.Range("I74").Value = 1446.0869
$ws.Range("J74").Value = 5690.6665
$ws.Range("K74").Value = 1446.0869
$ws.Range("L74").Value = 5690.6665
$ws.Range("M74").Value = -572.0869
$ws.Range("N74").Value = -7438.6665
$ws.Range("H77").Value = 2324.276
$ws.Range("I77").Value = 1446.0869
$ws.Range("J77").Value = 5690.6665
$ws.Range("K77").Value = 7230.4345
$ws.Range("L77").Value = 28453.3325
$ws.Range("M77").Value = -2862.4345
$ws.Range("N77").Value = -37189.3325
$ws.Range("H110").Value = 1652.8334
$ws.Range("I110").Value = 1257.6364
$ws.Range("K110").Value = 1257.6364
$ws.Range("M110").Value = 787.3635999999999
$ws.Range("H116").Value = 3400.8333
$ws.Range("I116").Value = 2882
$ws.Range("J116").Value = 5995
$ws.Range("K116").Value = 2882
$ws.Range("L116").Value = 5995
$ws.Range("M116").Value = -588
$ws.Range("N116").Value = -10583

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3400.8333
$ws.Range("I3").Value = 2882
$ws.Range("J3").Value = 5995
$ws.Range("K3").Value = 2882
$ws.Range("L3").Value = 5995
$ws.Range("M3").Value = -2768
$ws.Range("N3").Value = -6223
$ws.Range("H20").Value = 5969.75
$ws.Range("I20").Value = 5834.875
$ws.Range("J20").Value = 6104.625
$ws.Range("K20").Value = 5834.875
$ws.Range("L20").Value = 6104.625
$ws.Range("M20").Value = -5587.875
$ws.Range("N20").Value = -6598.625
$ws.Range("H94").Value = 56543.285
$ws.Range("I94").Value = 2025.8889
$ws.Range("K94").Value = 2025.8889
$ws.Range("M94").Value = -1574.8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 16496779
$ws.Range("I4").Value = 38726.5
$ws.Range("J4").Value = 23080000
$ws.Range("K4").Value = 38726.5
$ws.Range("L4").Value = 23080000
$ws.Range("M4").Value = -38614.5
$ws.Range("N4").Value = -23080224
$ws.Range("H31").Value = 2778.0527
$ws.Range("I31").Value = 1529.1852
$ws.Range("K31").Value = 1529.1852
$ws.Range("M31").Value = -1234.1852
$ws.Range("H34").Value = 2778.0527
$ws.Range("I34").Value = 1529.1852
$ws.Range("K34").Value = 1529.1852
$ws.Range("M34").Value = -1327.1852
$ws.Range("H62").Value = 3996.8
$ws.Range("I62").Value = 3996.8
$ws.Range("K62").Value = 3996.8
$ws.Range("M62").Value = -3372.8
$ws.Range("H65").Value = 3996.8
$ws.Range("I65").Value = 3996.8
$ws.Range("K65").Value = 19984
$ws.Range("M65").Value = -16864
$ws.Range("H86").Value = 58828230
$ws.Range("I86").Value = 100003670
$ws.Range("J86").Value = 6170
$ws.Range("K86").Value = 100003670
$ws.Range("L86").Value = 6170
$ws.Range("M86").Value = -100002547
$ws.Range("N86").Value = -8416
$ws.Range("H89").Value = 58828230
$ws.Range("I89").Value = 100003670
$ws.Range("J89").Value = 6170
$ws.Range("K89").Value = 500018350
$ws.Range("L89").Value = 30850
$ws.Range("M89").Value = -500012734
$ws.Range("N89").Value = -42082
$ws.Range("H105").Value = 681
$ws.Range("I105").Value = 703
$ws.Range("K105").Value = 703
$ws.Range("M105").Value = 1044
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 434.66666
$ws.Range("J23").Value = 464.8
$ws.Range("L23").Value = 1394.4
$ws.Range("N23").Value = -1864.4
$ws.Range("H93").Value = 147141.28
$ws.Range("I93").Value = 999999
$ws.Range("J93").Value = 4998.3335
$ws.Range("K93").Value = 2999997
$ws.Range("L93").Value = 14995.0005
$ws.Range("M93").Value = -2998125
$ws.Range("N93").Value = -18739.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 36078504
$ws.Range("I70").Value = 5555.3
$ws.Range("K70").Value = 5555.3
$ws.Range("M70").Value = -5285.3
$ws.Range("H73").Value = 36078504
$ws.Range("I73").Value = 5555.3
$ws.Range("K73").Value = 5555.3
$ws.Range("M73").Value = -4619.3
$ws.Range("H103").Value = 49797.5
$ws.Range("I103").Value = 49000
$ws.Range("J103").Value = 50595
$ws.Range("K103").Value = 49000
$ws.Range("L103").Value = 50595
$ws.Range("M103").Value = -47828
$ws.Range("N103").Value = -52939
$ws.Range("H122").Value = 1298.3334
$ws.Range("I122").Value = 1298.3334
$ws.Range("K122").Value = 3895.0002
$ws.Range("M122").Value = -1445.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2979
$ws.Range("I68").Value = 2723.75
$ws.Range("K68").Value = 2723.75
$ws.Range("M68").Value = -1974.75
$ws.Range("H71").Value = 2979
$ws.Range("I71").Value = 2723.75
$ws.Range("K71").Value = 13618.75
$ws.Range("M71").Value = -9874.75
$ws.Range("H93").Value = 31314.092
$ws.Range("I93").Value = 953
$ws.Range("K93").Value = 953
$ws.Range("M93").Value = 295
$ws.Range("H122").Value = 4154.3125
$ws.Range("I122").Value = 3890.3333
$ws.Range("J122").Value = 4946.25
$ws.Range("K122").Value = 11670.9999
$ws.Range("L122").Value = 14838.75
$ws.Range("M122").Value = -9220.999899999999
$ws.Range("N122").Value = -19738.75
$ws.Range("H132").Value = 4937.1333
$ws.Range("I132").Value = 4844.3184
$ws.Range("K132").Value = 14532.9552
$ws.Range("M132").Value = -12002.9552
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530
